# Updated cryptos list on Mon Mar 27 16:45:35 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and reorders two rows (WEMIXTOKEN now above TheSandbox).
# Price values are written with a leading apostrophe so Excel stores them
# as literal text (matching the source data, which keeps trailing zeros /
# dot-grouped thousands such as "1.000" or "27.189.39" that would otherwise
# be silently reinterpreted as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.189.39"
$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").Value = "'1.722.11"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'312.42"
$ws.Range("E5").Value = "  -4.57%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "'0.4610"
$ws.Range("E7").Value = "  +2.93%  "

$ws.Range("D8").Value = "'0.3447"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").Value = "'42.68"
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("D10").Value = "'0.07274"
$ws.Range("E10").Value = "  -2.23%  "

$ws.Range("D11").Value = "'1.046"
$ws.Range("E11").Value = "  -4.38%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'19.82"
$ws.Range("E13").Value = "  -4.84%  "

$ws.Range("D14").Value = "'5.855"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "'1.720.93"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").Value = "'6.879"
$ws.Range("E16").Value = "  -4.48%  "

$ws.Range("D17").Value = "'89.62"
$ws.Range("E17").Value = "  -3.50%  "

$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D19").Value = "'0.06313"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "'16.50"
$ws.Range("E21").Value = "  -3.68%  "

$ws.Range("D22").Value = "'5.624"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "'27.253.42"
$ws.Range("E23").Value = "  -2.31%  "

$ws.Range("D24").Value = "'10.88"
$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("D25").Value = "'2.130"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("D26").Value = "'154.92"
$ws.Range("E26").Value = "  -4.77%  "

$ws.Range("D27").Value = "'19.35"
$ws.Range("E27").Value = "  -4.21%  "

$ws.Range("D28").Value = "'1.921.20"
$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("D29").Value = "'2.165"
$ws.Range("E29").Value = "  -0.42%  "

$ws.Range("D30").Value = "'119.23"
$ws.Range("E30").Value = "  -4.66%  "

$ws.Range("D31").Value = "'1.034"
$ws.Range("E31").Value = "  -6.11%  "

$ws.Range("D32").Value = "'0.09104"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").Value = "'3.586"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").Value = "'5.352"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("D35").Value = "'0.02212"
$ws.Range("E35").Value = "  -3.49%  "

$ws.Range("D36").Value = "'0.05859"
$ws.Range("E36").Value = "  -4.06%  "

$ws.Range("D37").Value = "'11.11"
$ws.Range("E37").Value = "  -6.15%  "

$ws.Range("D38").Value = "'0.1999"
$ws.Range("E38").Value = "  -4.47%  "

$ws.Range("D39").Value = "'4.725"
$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "'1.402"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5959"
$ws.Range("E41").Value = "  -5.53%  "

$ws.Range("D42").Value = "'1.132"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("D43").Value = "'7.478"
$ws.Range("E43").Value = "  -5.80%  "

$ws.Range("D44").Value = "'12.76"
$ws.Range("E44").Value = "  -4.20%  "

$ws.Range("D45").Value = "'3.591"
$ws.Range("E45").Value = "  -3.84%  "

$ws.Range("D46").Value = "'0.5635"
$ws.Range("E46").Value = "  -3.78%  "

$ws.Range("D47").Value = "'119.68"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").Value = "'1.869"
$ws.Range("E48").Value = "  -4.18%  "

$ws.Range("D49").Value = "'0.06664"
$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("D50").Value = "'1.085"
$ws.Range("E50").Value = "  -4.48%  "

$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.02%  "
